# Auto-generated Excel COM-interop script applying the Phantom_Profits sheet update
# (recomputed currentAveragePrice / LevePrice / LeveProfit columns from the XIVAPI price refresh)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 350.6
$ws.Range("I4").Value = 337.75
$ws.Range("J4").Value = 402
$ws.Range("K4").Value = 337.75
$ws.Range("L4").Value = 402
$ws.Range("M4").Value = -223.75
$ws.Range("N4").Value = -630
$ws.Range("H38").Value = 1666789.5
$ws.Range("I38").Value = 2500054.5
$ws.Range("K38").Value = 7500163.5
$ws.Range("M38").Value = -7499791.5
$ws.Range("H51").Value = 7856.7144
$ws.Range("I51").Value = 8332.833000000001
$ws.Range("K51").Value = 8332.833000000001
$ws.Range("M51").Value = -7848.833000000001
$ws.Range("H70").Value = 5344.6
$ws.Range("I70").Value = 4142.857
$ws.Range("J70").Value = 6396.125
$ws.Range("K70").Value = 12428.571
$ws.Range("L70").Value = 19188.375
$ws.Range("M70").Value = -12158.571
$ws.Range("N70").Value = -19728.375
$ws.Range("H73").Value = 5344.6
$ws.Range("I73").Value = 4142.857
$ws.Range("J73").Value = 6396.125
$ws.Range("K73").Value = 12428.571
$ws.Range("L73").Value = 19188.375
$ws.Range("M73").Value = -11492.571
$ws.Range("N73").Value = -21060.375
$ws.Range("H80").Value = 7764
$ws.Range("I80").Value = 7394
$ws.Range("J80").Value = 7949
$ws.Range("K80").Value = 22182
$ws.Range("L80").Value = 23847
$ws.Range("M80").Value = -21184
$ws.Range("N80").Value = -25843
$ws.Range("H83").Value = 7764
$ws.Range("I83").Value = 7394
$ws.Range("J83").Value = 7949
$ws.Range("K83").Value = 66546
$ws.Range("L83").Value = 71541
$ws.Range("M83").Value = -61554
$ws.Range("N83").Value = -81525
$ws.Range("H86").Value = 4934.3335
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 4934.3335
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H100").Value = 2721
$ws.Range("I100").Value = 2608.923
$ws.Range("K100").Value = 2608.923
$ws.Range("M100").Value = -2067.923

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 701.8
$ws.Range("J4").Value = 399.5
$ws.Range("L4").Value = 399.5
$ws.Range("N4").Value = -631.5
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("N8").Value = 0
$ws.Range("H63").Value = 21179
$ws.Range("I63").Value = 21375.5
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 21375.5
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = -20689.5
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 21179
$ws.Range("I66").Value = 21375.5
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 106877.5
$ws.Range("L66").Value = 100000
$ws.Range("M66").Value = -103445.5
$ws.Range("N66").Value = -106864
$ws.Range("H97").Value = 968.8570999999999
$ws.Range("I97").Value = 213.66667
$ws.Range("K97").Value = 213.66667
$ws.Range("M97").Value = 282.33333
$ws.Range("H122").Value = 1391.0312
$ws.Range("I122").Value = 1444.1333
$ws.Range("K122").Value = 4332.3999
$ws.Range("M122").Value = -1882.3999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 594.25
$ws.Range("J94").Value = 726.2857
$ws.Range("L94").Value = 726.2857
$ws.Range("N94").Value = -1628.2857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2013.8889
$ws.Range("I31").Value = 2032.4286
$ws.Range("K31").Value = 2032.4286
$ws.Range("M31").Value = -1737.4286
$ws.Range("H34").Value = 2013.8889
$ws.Range("I34").Value = 2032.4286
$ws.Range("K34").Value = 2032.4286
$ws.Range("M34").Value = -1830.4286
$ws.Range("H120").Value = 45000
$ws.Range("J120").Value = 45000
$ws.Range("L120").Value = 45000
$ws.Range("N120").Value = -52258

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1443928
$ws.Range("I4").Value = 19832.6
$ws.Range("K4").Value = 59497.8
$ws.Range("M4").Value = -59385.8
$ws.Range("H9").Value = 4500
$ws.Range("I9").Value = 4000
$ws.Range("K9").Value = 12000
$ws.Range("M9").Value = -11776
$ws.Range("H10").Value = 154.08333
$ws.Range("I10").Value = 163.54546
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 490.63638
$ws.Range("L10").Value = 150
$ws.Range("M10").Value = -351.63638
$ws.Range("N10").Value = -428
$ws.Range("H11").Value = 607.1429000000001
$ws.Range("I11").Value = 607.1429000000001
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1821.4287
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1681.4287
$ws.Range("H12").Value = 387.9091
$ws.Range("I12").Value = 389.57144
$ws.Range("J12").Value = 385
$ws.Range("K12").Value = 1168.71432
$ws.Range("L12").Value = 1155
$ws.Range("M12").Value = -995.71432
$ws.Range("N12").Value = -1501
$ws.Range("H46").Value = 628.8333
$ws.Range("J46").Value = 644
$ws.Range("L46").Value = 1932
$ws.Range("N46").Value = -2114
$ws.Range("H113").Value = 2800.8462
$ws.Range("J113").Value = 2863.3
$ws.Range("L113").Value = 8589.900000000001
$ws.Range("N113").Value = -12929.9

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.2
$ws.Range("I2").Value = 96.888885
$ws.Range("K2").Value = 96.888885
$ws.Range("M2").Value = 16.111115
$ws.Range("H12").Value = 929857.4399999999
$ws.Range("I12").Value = 967250.5
$ws.Range("J12").Value = 880000
$ws.Range("K12").Value = 967250.5
$ws.Range("L12").Value = 880000
$ws.Range("M12").Value = -967110.5
$ws.Range("N12").Value = -880280
$ws.Range("H13").Value = 3650
$ws.Range("I13").Value = 1500
$ws.Range("K13").Value = 1500
$ws.Range("M13").Value = -1361
$ws.Range("H80").Value = 4333.3335
$ws.Range("I80").Value = 4333.3335
$ws.Range("K80").Value = 4333.3335
$ws.Range("M80").Value = -3335.3335
$ws.Range("H83").Value = 4333.3335
$ws.Range("I83").Value = 4333.3335
$ws.Range("K83").Value = 21666.6675
$ws.Range("M83").Value = -16674.6675
$ws.Range("H102").Value = 2023.4445
$ws.Range("I102").Value = 2026.375
$ws.Range("K102").Value = 2026.375
$ws.Range("M102").Value = -404.375
$ws.Range("H122").Value = 4230.857
$ws.Range("I122").Value = 3324.2
$ws.Range("J122").Value = 6497.5
$ws.Range("K122").Value = 9972.599999999999
$ws.Range("L122").Value = 19492.5
$ws.Range("M122").Value = -7522.599999999999
$ws.Range("N122").Value = -24392.5
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -6530

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6854
$ws.Range("I7").Value = 6854
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 6854
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -6742
$ws.Range("H61").Value = 1826
$ws.Range("I61").Value = 1826
$ws.Range("K61").Value = 1826
$ws.Range("M61").Value = -1624
$ws.Range("H100").Value = 845.625
$ws.Range("I100").Value = 796.6667
$ws.Range("J100").Value = 992.5
$ws.Range("K100").Value = 796.6667
$ws.Range("L100").Value = 992.5
$ws.Range("M100").Value = -255.6667
$ws.Range("N100").Value = -2074.5
$ws.Range("H113").Value = 1826
$ws.Range("I113").Value = 1826
$ws.Range("K113").Value = 1826
$ws.Range("M113").Value = 344
$ws.Range("H122").Value = 3923
$ws.Range("I122").Value = 3375.2778
$ws.Range("J122").Value = 5566.1665
$ws.Range("K122").Value = 10125.8334
$ws.Range("L122").Value = 16698.4995
$ws.Range("M122").Value = -7675.8334
$ws.Range("N122").Value = -21598.4995
$ws.Range("H126").Value = 6854
$ws.Range("I126").Value = 6854
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20562
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -18092
$ws.Range("H132").Value = 4891.143
$ws.Range("I132").Value = 4547.6
$ws.Range("K132").Value = 13642.8
$ws.Range("M132").Value = -11112.8
$ws.Range("H136").Value = 26318544
$ws.Range("I136").Value = 1474.5714
$ws.Range("J136").Value = 100006340
$ws.Range("K136").Value = 4423.7142
$ws.Range("L136").Value = 300019020
$ws.Range("M136").Value = -1873.7142
$ws.Range("N136").Value = -300024120

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17093
$ws.Range("I62").Value = 15100
$ws.Range("J62").Value = 17425.166
$ws.Range("K62").Value = 15100
$ws.Range("L62").Value = 17425.166
$ws.Range("M62").Value = -14476
$ws.Range("N62").Value = -18673.166
$ws.Range("H65").Value = 17093
$ws.Range("I65").Value = 15100
$ws.Range("J65").Value = 17425.166
$ws.Range("K65").Value = 75500
$ws.Range("L65").Value = 87125.83
$ws.Range("M65").Value = -72380
$ws.Range("N65").Value = -93365.83
$ws.Range("H96").Value = 1800
$ws.Range("I96").Value = 1800
$ws.Range("J96").Value = 1800
$ws.Range("K96").Value = 1800
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = -427
$ws.Range("N96").Value = -4546
$ws.Range("H126").Value = 3140.9
$ws.Range("I126").Value = 2900.75
$ws.Range("J126").Value = 4101.5
$ws.Range("K126").Value = 8702.25
$ws.Range("L126").Value = 12304.5
$ws.Range("M126").Value = -6232.25
$ws.Range("N126").Value = -17244.5
$ws.Range("H132").Value = 3199.88
$ws.Range("I132").Value = 2749.85
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8249.549999999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5719.549999999999
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 8580.857
$ws.Range("I136").Value = 6324.5
$ws.Range("K136").Value = 18973.5
$ws.Range("M136").Value = -16423.5

Write-Host "Phantom_Profits sheets updated."
